# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" on the
#   Overview sheet (both language columns) and on each language sheet's
#   Status column.
# - Refresh the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps
#   to the new handoff run.
# - Widen the Status columns so the longer "Ready for handoff" label fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status + latest generate date ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-27 16:57:05"

# --- zh-cn detail sheet: Status + Latest Handoff Datetime ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-27 16:56:58"

# --- de-de detail sheet: Status + Latest Handoff Datetime ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-27 16:57:05"

# --- Widen the Status columns to fit "Ready for handoff" ---
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
